$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crendetials")

# Row 7 used to hold the "Server" credentials; replace it with the Putty
# credentials that used to live in row 11 (which is being removed below).
$ws.Range("C7").Value = "Putty"
$ws.Range("D7").Value = "132.148.72.192"
$ws.Range("E7").Value = "esc"
$ws.Range("F7").Value = "Esc@esc123"

# Row 8 (phpmyadmin) gets a new host/URL and the password is now "root".
$ws.Range("D8").Value = "http://132.148.72.192/phpmyadmin/"
$ws.Range("F8").Value = "root"

# The old Putty row (row 11) is no longer needed, its data moved to row 7.
$ws.Rows.Item(11).Delete()

# Remove the now-unused extra worksheets.
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Restore the last-known selection.
$ws.Range("D17").Select()
